$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2025-05-05 Monday", $true, $false, $false, $false, $false, $true, 1, $false, "2025-05-06 Tuesday", 2) | Out-Null

# Update each math-problem cell in the table, by position (row, column)
$t = $d.Tables.Item(1)

$t.Rows.Item(1).Cells.Item(1).Range.Text = "18+14="
$t.Rows.Item(1).Cells.Item(2).Range.Text = "53+38="
$t.Rows.Item(1).Cells.Item(3).Range.Text = "94-16="
$t.Rows.Item(1).Cells.Item(4).Range.Text = "22+69="
$t.Rows.Item(1).Cells.Item(5).Range.Text = "96-48="
$t.Rows.Item(2).Cells.Item(1).Range.Text = "81-42="
$t.Rows.Item(2).Cells.Item(2).Range.Text = "73-25="
$t.Rows.Item(2).Cells.Item(3).Range.Text = "47+27="
$t.Rows.Item(2).Cells.Item(4).Range.Text = "76-38="
$t.Rows.Item(2).Cells.Item(5).Range.Text = "68+25="
$t.Rows.Item(3).Cells.Item(1).Range.Text = "77-8="
$t.Rows.Item(3).Cells.Item(2).Range.Text = "9+6="
$t.Rows.Item(3).Cells.Item(3).Range.Text = "14+8="
$t.Rows.Item(3).Cells.Item(4).Range.Text = "52-47="
$t.Rows.Item(3).Cells.Item(5).Range.Text = "74-8="
$t.Rows.Item(4).Cells.Item(1).Range.Text = "68+7="
$t.Rows.Item(4).Cells.Item(2).Range.Text = "84-18="
$t.Rows.Item(4).Cells.Item(3).Range.Text = "17+7="
$t.Rows.Item(4).Cells.Item(4).Range.Text = "41-24="
$t.Rows.Item(4).Cells.Item(5).Range.Text = "18+5="
$t.Rows.Item(5).Cells.Item(1).Range.Text = "73-9="
$t.Rows.Item(5).Cells.Item(2).Range.Text = "20-13="
$t.Rows.Item(5).Cells.Item(3).Range.Text = "3+58="
$t.Rows.Item(5).Cells.Item(4).Range.Text = "18+25="
$t.Rows.Item(5).Cells.Item(5).Range.Text = "9+69="
$t.Rows.Item(6).Cells.Item(1).Range.Text = "59+29="
$t.Rows.Item(6).Cells.Item(2).Range.Text = "24-18="
$t.Rows.Item(6).Cells.Item(3).Range.Text = "75-58="
$t.Rows.Item(6).Cells.Item(4).Range.Text = "87-58="
$t.Rows.Item(6).Cells.Item(5).Range.Text = "21-9="
$t.Rows.Item(7).Cells.Item(1).Range.Text = "53-17="
$t.Rows.Item(7).Cells.Item(2).Range.Text = "17+77="
$t.Rows.Item(7).Cells.Item(3).Range.Text = "71-53="
$t.Rows.Item(7).Cells.Item(4).Range.Text = "3+48="
$t.Rows.Item(7).Cells.Item(5).Range.Text = "85+6="
$t.Rows.Item(8).Cells.Item(1).Range.Text = "17+17="
$t.Rows.Item(8).Cells.Item(2).Range.Text = "70-17="
$t.Rows.Item(8).Cells.Item(3).Range.Text = "60-11="
$t.Rows.Item(8).Cells.Item(4).Range.Text = "27+4="
$t.Rows.Item(8).Cells.Item(5).Range.Text = "86-57="
$t.Rows.Item(9).Cells.Item(1).Range.Text = "81-77="
$t.Rows.Item(9).Cells.Item(2).Range.Text = "65+8="
$t.Rows.Item(9).Cells.Item(3).Range.Text = "76+7="
$t.Rows.Item(9).Cells.Item(4).Range.Text = "58+34="
$t.Rows.Item(9).Cells.Item(5).Range.Text = "47+48="
$t.Rows.Item(10).Cells.Item(1).Range.Text = "34+39="
$t.Rows.Item(10).Cells.Item(2).Range.Text = "50-23="
$t.Rows.Item(10).Cells.Item(3).Range.Text = "74-26="
$t.Rows.Item(10).Cells.Item(4).Range.Text = "90-72="
$t.Rows.Item(10).Cells.Item(5).Range.Text = "28+53="
$t.Rows.Item(11).Cells.Item(1).Range.Text = "6+36="
$t.Rows.Item(11).Cells.Item(2).Range.Text = "16+38="
$t.Rows.Item(11).Cells.Item(3).Range.Text = "98-59="
$t.Rows.Item(11).Cells.Item(4).Range.Text = "96-59="
$t.Rows.Item(11).Cells.Item(5).Range.Text = "53-28="
$t.Rows.Item(12).Cells.Item(1).Range.Text = "36-29="
$t.Rows.Item(12).Cells.Item(2).Range.Text = "40-33="
$t.Rows.Item(12).Cells.Item(3).Range.Text = "73-57="
$t.Rows.Item(12).Cells.Item(4).Range.Text = "29+7="
$t.Rows.Item(12).Cells.Item(5).Range.Text = "69+22="
$t.Rows.Item(13).Cells.Item(1).Range.Text = "80-7="
$t.Rows.Item(13).Cells.Item(2).Range.Text = "42-18="
$t.Rows.Item(13).Cells.Item(3).Range.Text = "5+76="
$t.Rows.Item(13).Cells.Item(4).Range.Text = "72+9="
$t.Rows.Item(13).Cells.Item(5).Range.Text = "80-24="
$t.Rows.Item(14).Cells.Item(1).Range.Text = "3+48="
$t.Rows.Item(14).Cells.Item(2).Range.Text = "33+28="
$t.Rows.Item(14).Cells.Item(3).Range.Text = "35+19="
$t.Rows.Item(14).Cells.Item(4).Range.Text = "7+79="
$t.Rows.Item(14).Cells.Item(5).Range.Text = "69+18="
$t.Rows.Item(15).Cells.Item(1).Range.Text = "79+13="
$t.Rows.Item(15).Cells.Item(2).Range.Text = "91-69="
$t.Rows.Item(15).Cells.Item(3).Range.Text = "62-37="
$t.Rows.Item(15).Cells.Item(4).Range.Text = "26+15="
$t.Rows.Item(15).Cells.Item(5).Range.Text = "82-6="
$t.Rows.Item(16).Cells.Item(1).Range.Text = "95-29="
$t.Rows.Item(16).Cells.Item(2).Range.Text = "8+17="
$t.Rows.Item(16).Cells.Item(3).Range.Text = "50-39="
$t.Rows.Item(16).Cells.Item(4).Range.Text = "5+49="
$t.Rows.Item(16).Cells.Item(5).Range.Text = "81-9="
$t.Rows.Item(17).Cells.Item(1).Range.Text = "43-28="
$t.Rows.Item(17).Cells.Item(2).Range.Text = "28-19="
$t.Rows.Item(17).Cells.Item(3).Range.Text = "20-4="
$t.Rows.Item(17).Cells.Item(4).Range.Text = "16+67="
$t.Rows.Item(17).Cells.Item(5).Range.Text = "47+47="
$t.Rows.Item(18).Cells.Item(1).Range.Text = "36+27="
$t.Rows.Item(18).Cells.Item(2).Range.Text = "58+38="
$t.Rows.Item(18).Cells.Item(3).Range.Text = "88+9="
$t.Rows.Item(18).Cells.Item(4).Range.Text = "40-14="
$t.Rows.Item(18).Cells.Item(5).Range.Text = "46+35="
$t.Rows.Item(19).Cells.Item(1).Range.Text = "58+6="
$t.Rows.Item(19).Cells.Item(2).Range.Text = "28+58="
$t.Rows.Item(19).Cells.Item(3).Range.Text = "79+15="
$t.Rows.Item(19).Cells.Item(4).Range.Text = "84-45="
$t.Rows.Item(19).Cells.Item(5).Range.Text = "71-19="
$t.Rows.Item(20).Cells.Item(1).Range.Text = "81-57="
$t.Rows.Item(20).Cells.Item(2).Range.Text = "45+29="
$t.Rows.Item(20).Cells.Item(3).Range.Text = "18+74="
$t.Rows.Item(20).Cells.Item(4).Range.Text = "5+29="
$t.Rows.Item(20).Cells.Item(5).Range.Text = "82-76="
